$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.376.70"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.89%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.537.34"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.19%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.77%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.03"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.567"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.27%  "
$ws.Range("E8").Value = "  +0.25%  "
$ws.Range("E9").Value = "  -3.64%  "
$ws.Range("E10").Value = "  -0.74%  "
$ws.Range("E11").Value = "  -1.37%  "
$ws.Range("E12").Value = "  -2.00%  "
$ws.Range("E13").Value = "  -0.20%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.928.10"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.03%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.94"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.55%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.503.13"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.00%  "
$ws.Range("E17").Value = "  -1.61%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.383.96"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.76%  "
$ws.Range("E19").Value = "  -1.33%  "
$ws.Range("E20").Value = "  -1.85%  "
$ws.Range("E21").Value = "  -3.92%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.76"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.88%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "242.76"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.67%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.90"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.79%  "
$ws.Range("E25").Value = "  -1.14%  "
$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.36"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.62%  "
$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.06%  "
$ws.Range("E28").Value = "  -3.91%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "39.76"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.62%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.13"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.68%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "157.60"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.69%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.72"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.99%  "
$ws.Range("E33").Value = "  +16.85%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0796"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.30%  "
$ws.Range("E35").Value = "  -2.91%  "
$ws.Range("E36").Value = "  -5.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.16"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.21%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.02"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.53%  "
$ws.Range("E39").Value = "  -1.67%  "
$ws.Range("E40").Value = "  -0.86%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.21"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +8.48%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.78"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.01"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.44%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.29"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.94%  "
$ws.Range("E45").Value = "  -2.96%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.954.91"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.63%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.91"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.39%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.782.11"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.79%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "80.73"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.83%  "
$ws.Range("E50").Value = "  -1.38%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "101.27"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.19%  "
